$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Start Time on row 7 (was "10/25/2021 05:24 P.M.", now an actual
# start time for the new gantt-chart work session).
$ws.Range("A7").Value = "10/27/2021 1:00 P.M."

# Fill in the End Time for row 7.
$ws.Range("B7").Value = "10/27/2021 2:00 P.M."

# Time spent on row 7: 1:00:00 (1 hour).
$ws.Range("C7").Value = (1*3600 + 0*60 + 0)/3600/24

# Running total time through row 7: 5:10:00.
$ws.Range("D7").Value = (5*3600 + 10*60 + 0)/3600/24

# Summary of work for row 7.
$ws.Range("E7").Value = "Created gantt chart and made additional small changes to scope"

# Correct the running total on row 6 to 4:10:00.
$ws.Range("D6").Value = (4*3600 + 10*60 + 0)/3600/24

# Leave the final selection on the newly filled-in summary cell.
$ws.Range("E7").Select() | Out-Null
